$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgabenliste")

# Mark rows 17 and 18 as "done" (B column uses the "a" marker glyph already
# used elsewhere in the sheet) and add remark text in column E for the two
# new tasks (Meth1 / MethAttr).
$ws.Range("B17").Value = "a"
$ws.Range("E17").Value = "Meth1"

$ws.Range("B18").Value = "a"
$ws.Range("E18").Value = "MethAttr"

# Update the active sheet's selection to match the author's saved view.
$ws.Activate()
$ws.Range("C19").Select()
